$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the Wing Area calculation formula in D2
$ws.Range("D2").Formula = "=(B14*B15)+(((B17-B15)*B14)*0.5)"

# Remove the old "Center of Gravity x coord" row (row 4) entirely
$ws.Range("A4:C4").ClearContents()

# Update the view: scroll back to top and move the selection to G8
$null = $ws.Range("G8").Select()
